# Weekly data refresh: a new price record for
# "Terminal La Palmera de La Serena - Zapallo italiano" is inserted as the
# new row 159, pushing the previously-existing rows 159-190 down to 160-191
# (the used range grows from A1:R190 to A1:R191).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 159; Excel shifts rows 159-190 down to
# 160-191 and extends the sheet's used range accordingly.
$ws.Rows(159).Insert()

# Fill in the newly inserted row 159 with the new weekly record.
$ws.Cells.Item(159, 1).Value  = 8
$ws.Cells.Item(159, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(159, 3).Value  = "Coquimbo"
$ws.Cells.Item(159, 4).Value  = 44476
$ws.Cells.Item(159, 5).Value  = 4
$ws.Cells.Item(159, 6).Value  = 100112032
$ws.Cells.Item(159, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(159, 8).Value  = "Sin especificar"
$ws.Cells.Item(159, 9).Value  = "Primera"
$ws.Cells.Item(159, 10).Value = 520
$ws.Cells.Item(159, 11).Value = 19000
$ws.Cells.Item(159, 12).Value = 20000
$ws.Cells.Item(159, 13).Value = 19500
$ws.Cells.Item(159, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(159, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(159, 16).Value = 279
$ws.Cells.Item(159, 17).Value = 70
$ws.Cells.Item(159, 18).Value = "Hortaliza"
